$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 54) below the existing data, matching the
# formatting of the last existing row (row 53) in column A.
$ws.Range("A53").Copy()
$ws.Range("A54").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(54, 1).Value = 45986
$ws.Cells.Item(54, 2).Value = 2025
$ws.Cells.Item(54, 3).Value = 2.043309689777173
$ws.Cells.Item(54, 4).Value = 2026
$ws.Cells.Item(54, 5).Value = 0.9040423720836799
